# Weekly update: insert a new week's price record for
# "Terminal La Palmera de La Serena - Zanahoria" at row 274, pushing the
# existing rows 274:290 down to 275:291 (last row becomes row 291).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 274 - existing rows 274..290 shift to 275..291
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row 274 with this week's data
$ws.Cells.Item(274, 1).Value  = 8
$ws.Cells.Item(274, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(274, 3).Value  = "Coquimbo"
$ws.Cells.Item(274, 4).Value  = 44610
$ws.Cells.Item(274, 5).Value  = 4
$ws.Cells.Item(274, 6).Value  = 100114013
$ws.Cells.Item(274, 7).Value  = "Zanahoria"
$ws.Cells.Item(274, 8).Value  = "Sin especificar"
$ws.Cells.Item(274, 9).Value  = "Primera"
$ws.Cells.Item(274, 10).Value = 600
$ws.Cells.Item(274, 11).Value = 5500
$ws.Cells.Item(274, 12).Value = 6000
$ws.Cells.Item(274, 13).Value = 5750
$ws.Cells.Item(274, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(274, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(274, 16).Value = 288
$ws.Cells.Item(274, 17).Value = 20
$ws.Cells.Item(274, 18).Value = "Hortaliza"

# Copy the date cell's number format from the row below (preserves the
# date display format for the shifted-down rows / new row).
$ws.Cells.Item(274, 4).NumberFormat = $ws.Cells.Item(275, 4).NumberFormat
